$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A5 value (existing record correction)
$ws.Range("A5").Value = 171028

# Add new student rows
$ws.Range("A7").Value = 172064
$ws.Range("B7").Value = "Aleksandar"
$ws.Range("C7").Value = "Velickovski"

$ws.Range("A8").Value = 175032
$ws.Range("B8").Value = "Predrag"
$ws.Range("C8").Value = "Spasovski"

# Update selection to match author's final cursor position
$ws.Range("A5").Select()
